$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: fill in ДЗ_6 (H2) and a plain 7 value (I2) ---
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 7

# --- Row 24: fill in G24, H24 and new I24 with the green/thick-sided style ---
$ws.Range("G24").Value = 5
$ws.Range("H24").Value = 5

$i24 = $ws.Range("I24")
$i24.Value = 5
$i24.Interior.Color = 5296274
$i24.HorizontalAlignment = -4108
$i24.VerticalAlignment = -4108
$i24.WrapText = $true
$i24.Borders.Item(7).Weight = 4
$i24.Borders.Item(7).Color = 0
$i24.Borders.Item(10).Weight = 4
$i24.Borders.Item(10).Color = 0

# --- J24: turn the shared total formula into an explicit one that also covers I24 ---
$ws.Range("J24").Formula = "=SUM(C24:I24)"

# --- Update the active selection to J6 ---
$ws.Range("J6").Select()
